$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row labels: _old -> _FV2404, _new -> _FV2410
for ($i = 0; $i -lt 10; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2404")
}
for ($i = 0; $i -lt 10; $i++) {
    $cell = $ws.Cells.Item(1, $i + 12)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2410")
}

# Freeze header row (pane split after row 1, frozen)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a table
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U94"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
